$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly garlic price record was added to the data set. In the
# source data (ordered by date) this new record belongs right before the
# current row 277, so insert a fresh row there; Excel will push the
# existing rows 277-298 down to 278-299 and extend the used range
# accordingly.
$dateFormat = $ws.Cells.Item(277, 4).NumberFormat
$ws.Rows.Item(277).Insert()

$ws.Cells.Item(277, 1).Value = 11
$ws.Cells.Item(277, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(277, 3).Value = "Bíobío"
$ws.Cells.Item(277, 4).Value2 = 45106
$ws.Cells.Item(277, 5).Value = 8
$ws.Cells.Item(277, 6).Value = 100112003
$ws.Cells.Item(277, 7).Value = "Ajo"
$ws.Cells.Item(277, 8).Value = "Chino"
$ws.Cells.Item(277, 9).Value = "Primera"
$ws.Cells.Item(277, 10).Value = 200
$ws.Cells.Item(277, 11).Value = 16000
$ws.Cells.Item(277, 12).Value = 17000
$ws.Cells.Item(277, 13).Value = 16500
$ws.Cells.Item(277, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(277, 15).Value = "China"
$ws.Cells.Item(277, 16).Value = 1650
$ws.Cells.Item(277, 17).Value = 10
$ws.Cells.Item(277, 18).Value = "Hortaliza"

# Make sure the new row's date cell uses the same date/time display
# format as the rest of column D.
$ws.Cells.Item(277, 4).NumberFormat = $dateFormat
